$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a text value, preventing Excel from auto-converting
# numeric-looking strings (e.g. "0.999") into floating point numbers, while
# keeping the cell style identical to its original (no quote-prefix formatting).
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for rows with both changes
Set-TextValue $ws.Range("D2") "90.927.78"
$ws.Range("E2").Value = "  +3.71%  "
Set-TextValue $ws.Range("D3") "3.171.42"
$ws.Range("E3").Value = "  -0.25%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.41%  "
Set-TextValue $ws.Range("D5") "215.97"
$ws.Range("E5").Value = "  +3.78%  "
Set-TextValue $ws.Range("D6") "629.21"
$ws.Range("E6").Value = "  +3.18%  "
Set-TextValue $ws.Range("D7") "0.396"
$ws.Range("E7").Value = "  +2.07%  "
Set-TextValue $ws.Range("D8") "0.713"
$ws.Range("E8").Value = "  +6.09%  "
Set-TextValue $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  -0.14%  "
Set-TextValue $ws.Range("D10") "3.164.14"
$ws.Range("E10").Value = "  -0.53%  "
Set-TextValue $ws.Range("D11") "0.565"
$ws.Range("E11").Value = "  +4.48%  "
Set-TextValue $ws.Range("D12") "0.180"
$ws.Range("E12").Value = "  +2.06%  "
Set-TextValue $ws.Range("D13") "0.0000254"
$ws.Range("E13").Value = "  +3.77%  "
Set-TextValue $ws.Range("D14") "90.580.04"
$ws.Range("E14").Value = "  +3.33%  "
Set-TextValue $ws.Range("D15") "5.31"
$ws.Range("E15").Value = "  +0.47%  "
Set-TextValue $ws.Range("D16") "3.750.22"
$ws.Range("E16").Value = "  -0.57%  "
Set-TextValue $ws.Range("D17") "32.56"
$ws.Range("E17").Value = "  +0.49%  "
Set-TextValue $ws.Range("D18") "3.173.10"
$ws.Range("E18").Value = "  -0.33%  "
Set-TextValue $ws.Range("D19") "3.30"
$ws.Range("E19").Value = "  +1.75%  "
Set-TextValue $ws.Range("D20") "0.0000213"
$ws.Range("E20").Value = "  +60.27%  "
Set-TextValue $ws.Range("D21") "439.43"
$ws.Range("E21").Value = "  +6.26%  "
Set-TextValue $ws.Range("D22") "13.28"
$ws.Range("E22").Value = "  -1.71%  "
Set-TextValue $ws.Range("D23") "8.46"
$ws.Range("E23").Value = "  -0.67%  "
Set-TextValue $ws.Range("D24") "4.98"
$ws.Range("E24").Value = "  -2.11%  "
Set-TextValue $ws.Range("D25") "5.19"
$ws.Range("E25").Value = "  -1.52%  "
Set-TextValue $ws.Range("D26") "11.64"
$ws.Range("E26").Value = "  -5.20%  "
Set-TextValue $ws.Range("D27") "80.05"
$ws.Range("E27").Value = "  +8.79%  "
Set-TextValue $ws.Range("D28") "3.336.58"
$ws.Range("E28").Value = "  -0.37%  "
Set-TextValue $ws.Range("D31") "0.156"
$ws.Range("E31").Value = "  -4.34%  "
Set-TextValue $ws.Range("D32") "4.03"
$ws.Range("E32").Value = "  +31.80%  "
Set-TextValue $ws.Range("D33") "8.28"
$ws.Range("E33").Value = "  +0.38%  "
Set-TextValue $ws.Range("D34") "520.71"
$ws.Range("E34").Value = "  -4.62%  "
Set-TextValue $ws.Range("D35") "6.93"
$ws.Range("E35").Value = "  +0.17%  "
Set-TextValue $ws.Range("D36") "1.88"
$ws.Range("E36").Value = "  +1.01%  "
Set-TextValue $ws.Range("D37") "1.29"
$ws.Range("E37").Value = "  -2.50%  "
Set-TextValue $ws.Range("D38") "22.33"
$ws.Range("E38").Value = "  +2.02%  "
Set-TextValue $ws.Range("D40") "0.998"
$ws.Range("E40").Value = "  -0.30%  "
Set-TextValue $ws.Range("D41") "0.126"
$ws.Range("E41").Value = "  -4.45%  "
Set-TextValue $ws.Range("D43") "1.91"
$ws.Range("E43").Value = "  -0.45%  "
Set-TextValue $ws.Range("D44") "0.368"
$ws.Range("E44").Value = "  -1.59%  "
Set-TextValue $ws.Range("D45") "146.64"
$ws.Range("E45").Value = "  -2.31%  "
Set-TextValue $ws.Range("D46") "44.05"
$ws.Range("E46").Value = "  +1.72%  "
Set-TextValue $ws.Range("D47") "170.83"
$ws.Range("E47").Value = "  -1.99%  "
Set-TextValue $ws.Range("D48") "0.126"
$ws.Range("E48").Value = "  +0.62%  "
Set-TextValue $ws.Range("D49") "0.739"
$ws.Range("E49").Value = "  +6.24%  "
Set-TextValue $ws.Range("D50") "24.62"
$ws.Range("E50").Value = "  +2.73%  "

# Update Volume(1h) (E) column only for rows without a price change
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("E42").Value = "  -0.08%  "

# Row 51: coin changed from ImmutableX to ARBITRUM (name, link, price, volume)
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D51") "0.609"
$ws.Range("E51").Value = "  +2.49%  "
